$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title 1: "A" + " " + "slide" -> single run "A slide"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "A slide"

# Content Placeholder 5 (table): cell (1,2) "a" + " " + "table" -> single run "a table"
$s.Shapes.Item(3).Table.Cell(1, 2).Shape.TextFrame.TextRange.Text = "a table"

# TextBox 3: "Plus" + " " + "an" + " " + "image" -> single run "Plus an image"
$s.Shapes.Item(7).TextFrame.TextRange.Text = "Plus an image"
